# Apply the edit described by the commit "two testcases 9 10"
# 1. Fill column B of sheet "Pruthvi" with catalog URLs (rows 1-9).
# 2. Add a new worksheet "Pruthvi1" with the first two URLs in column A.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$urls = @(
    "/Manufacturing-Processing-Machinery-Catalog/Machine-Tools.html",
    "/Manufacturing-Processing-Machinery-Catalog/Engineering-Construction-Machinery.html",
    "/Manufacturing-Processing-Machinery-Catalog/Woodworking-Machinery.html",
    "/Manufacturing-Processing-Machinery-Catalog/Plastic-Machinery.html",
    "/Manufacturing-Processing-Machinery-Catalog/Metallic-Processing-Machinery.html",
    "/Manufacturing-Processing-Machinery-Catalog/Mould.html",
    "/Manufacturing-Processing-Machinery-Catalog/Laser-Equipment.html",
    "/Manufacturing-Processing-Machinery-Catalog/Casting-Forging.html",
    "/Manufacturing-Processing-Machinery-Catalog/Agricultural-Machinery.html"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 1
    $ws1.Cells.Item($row, 2).Value = $urls[$i]
}

# Add the new worksheet after the existing one and rename it.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Pruthvi1"

$ws2.Cells.Item(1, 1).Value = $urls[0]
$ws2.Cells.Item(2, 1).Value = $urls[1]
